$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column keeps being stored as text (values like "289.68"
# or "1.544.96" would otherwise be auto-converted to numbers by Excel).
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range("D2").Value = "21.719.71"
$ws.Range("E2").Value = "  -1.31%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.539.77"
$ws.Range("E3").Value = "  -0.90%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.02%  "

# Row 5 - USDC
$ws.Range("E5").Value = "  +0.07%  "

# Row 6 - BNB
$ws.Range("D6").Value = "289.68"
$ws.Range("E6").Value = "  +1.11%  "

# Row 7 - XRP
$ws.Range("D7").Value = "0.3904"
$ws.Range("E7").Value = "  +3.45%  "

# Row 8 - Cardano
$ws.Range("D8").Value = "0.3167"

# Row 9 - OKB
$ws.Range("D9").Value = "42.90"
$ws.Range("E9").Value = "  +4.20%  "

# Row 10 - Dogecoin
$ws.Range("D10").Value = "0.07171"
$ws.Range("E10").Value = "  -1.74%  "

# Row 11 - Polygon
$ws.Range("E11").Value = "  -6.03%  "

# Row 12 - BinanceUSD
$ws.Range("E12").Value = "  +0.03%  "

# Row 13 - Polkadot
$ws.Range("D13").Value = "5.612"
$ws.Range("E13").Value = "  -1.72%  "

# Row 14 - Solana
$ws.Range("D14").Value = "18.56"
$ws.Range("E14").Value = "  -3.96%  "

# Row 15 - Chainlink
$ws.Range("D15").Value = "6.611"
$ws.Range("E15").Value = "  -2.77%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "1.544.96"
$ws.Range("E16").Value = "  -0.63%  "

# Row 17 - ShibaInu
$ws.Range("D17").Value = "0.00001098"
$ws.Range("E17").Value = "  +1.73%  "

# Row 18 - TRON
$ws.Range("D18").Value = "0.06568"
$ws.Range("E18").Value = "  -1.24%  "

# Row 19 - Litecoin
$ws.Range("D19").Value = "83.03"
$ws.Range("E19").Value = "  -2.35%  "

# Row 20 - Dai
$ws.Range("D20").Value = "1.001"
$ws.Range("E20").Value = "  +0.05%  "

# Row 21 - Uniswap
$ws.Range("D21").Value = "6.141"
$ws.Range("E21").Value = "  -4.36%  "

# Row 22 - Avalanche
$ws.Range("D22").Value = "15.33"
$ws.Range("E22").Value = "  -3.95%  "

# Row 23 - Cosmos
$ws.Range("E23").Value = "  -5.61%  "

# Row 24 - Toncoin
$ws.Range("D24").Value = "2.414"
$ws.Range("E24").Value = "  +6.26%  "

# Row 25 - WrappedBTC
$ws.Range("D25").Value = "21.723.20"
$ws.Range("E25").Value = "  -1.39%  "

# Row 26 - LidoDAOToken
$ws.Range("D26").Value = "2.356"
$ws.Range("E26").Value = "  -6.11%  "

# Row 27 - Monero
$ws.Range("D27").Value = "147.00"
$ws.Range("E27").Value = "  -1.91%  "

# Row 28 - EthereumClassic
$ws.Range("D28").Value = "18.36"
$ws.Range("E28").Value = "  -2.41%  "

# Row 29 - HuobiToken
$ws.Range("D29").Value = "4.848"
$ws.Range("E29").Value = "  +0.20%  "

# Row 30 - WrappedliquidstakedEther2.0
$ws.Range("D30").Value = "1.721.65"
$ws.Range("E30").Value = "  -0.42%  "

# Row 31 - BitcoinCash
$ws.Range("D31").Value = "117.43"
$ws.Range("E31").Value = "  -2.21%  "

# Row 32 & 33 - swap Filecoin and ImmutableX (rank numbers A32/A33 stay the same)
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "0.9631"
$ws.Range("E32").Value = "  -14.02%  "

$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "5.880"
$ws.Range("E33").Value = "  -0.50%  "

# Row 34 - Stellar
$ws.Range("D34").Value = "0.08171"
$ws.Range("E34").Value = "  +0.04%  "

# Row 35 - FraxShare
$ws.Range("D35").Value = "8.792"
$ws.Range("E35").Value = "  -5.23%  "

# Row 36 - Hedera
$ws.Range("D36").Value = "0.06060"
$ws.Range("E36").Value = "  -1.49%  "

# Row 37 - InternetComputer(DFINITY)
$ws.Range("D37").Value = "5.099"
$ws.Range("E37").Value = "  -2.45%  "

# Row 38 - VeChain
$ws.Range("D38").Value = "0.02196"
$ws.Range("E38").Value = "  -3.80%  "

# Row 39 - Algorand
$ws.Range("D39").Value = "0.2033"
$ws.Range("E39").Value = "  -3.57%  "

# Row 40 & 41 - swap WEMIXTOKEN and TrustWalletToken
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").Value = "1.176"
$ws.Range("E40").Value = "  -3.21%  "

$ws.Range("B41").Value = "WEMIXTOKEN"
$ws.Range("C41").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D41").Value = "1.425"
$ws.Range("E41").Value = "  -13.48%  "

# Row 42 - Frax
$ws.Range("D42").Value = "1.001"
$ws.Range("E42").Value = "  +0.04%  "

# Row 43 - Aptos
$ws.Range("E43").Value = "  -2.32%  "

# Row 44 - TheSandbox
$ws.Range("D44").Value = "0.5710"
$ws.Range("E44").Value = "  -3.61%  "

# Row 45 - PancakeSwap
$ws.Range("D45").Value = "3.735"
$ws.Range("E45").Value = "  +0.34%  "

# Row 46 - EnergySwap
$ws.Range("D46").Value = "12.90"
$ws.Range("E46").Value = "  -4.17%  "

# Row 47 - Decentraland
$ws.Range("E47").Value = "  -4.50%  "

# Row 48 - EOS
$ws.Range("D48").Value = "1.160"
$ws.Range("E48").Value = "  +0.63%  "

# Row 49 - Quant
$ws.Range("D49").Value = "115.97"
$ws.Range("E49").Value = "  -3.22%  "

# Row 50 - NEARProtocol
$ws.Range("D50").Value = "1.860"
$ws.Range("E50").Value = "  -3.64%  "

# Row 51 - Cronos
$ws.Range("E51").Value = "  -2.88%  "
